$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks first so Insert() does not leave stale/duplicated
# hyperlink relationships tied to the wrong cells once rows shift down.
$ws.Hyperlinks.Delete()

# Insert a new blank row at row 2, pushing the existing data rows (old 2-9)
# down to rows 3-10.
$ws.Rows.Item(2).Insert()

# Copy the formatting (styles/number formats) of the row below (the old row
# 2, now row 3) into the freshly inserted row 2 so it matches the rest of
# the data rows.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate new row 2 with the latest price entry.
$ws.Cells.Item(2, 1).Value = 9
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 274.95
$ws.Cells.Item(2, 5).Value = "30-09-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf"

# Re-create all hyperlinks (rows 2-10) so their relationships/targets match
# the cell text after the shift.
$urls = @(
    "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf",
    "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $url = $urls[$i]
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $url, "", "", $url) | Out-Null
}

# Hyperlinks.Add() stamps its cells with Excel's built-in blue/underlined
# "Hyperlink" style. The source file keeps the plain data-row style on the
# Circular Link column, so restore it across the whole column.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("F2:F10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
